$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 already carries the exact formatting pattern (border, no-wrap for
# A/B/D/E, wrap for C) that the two new test-case rows need, so clone it
# down into rows 29 and 30 before writing values.
$ws.Range("A10:E10").Copy()
$ws.Range("A29:E29").PasteSpecial(-4122)

$ws.Range("A10:E10").Copy()
$ws.Range("A30:E30").PasteSpecial(-4122)

# Values are entered Jira-id/Description first (for both rows), then the
# TCID column last, matching the shared-string insertion order recorded in
# the edited workbook.
$ws.Range("B29").Value = "OPQA-5154||OPQA-5230"
$ws.Range("C29").Value = "Verify that a ""Call us"" section is present in customer care contact page with customer care contact details (region, phone numbers, hours of operation, language supported||Ensure that the page has ""Support Request"" and ""Call us"" sections."

$ws.Range("B30").Value = "OPQA-5168 || OPQA-5227"
$webFormText = "Verify that the web form provided to user should be application specific with following required fields`n1.Name 2.Organization 3.Contact details (email, telephone) 4.Issue Category 5.Country`n6.Description of issue ( a free form text box where a user can describe why they are contacting support) ||`nVerify that the web form provided to user will be application specific with following required fields (Name, Organization, email, telephone, Issue Category, Country, Description of issue)"
$ws.Range("C30").Value = $webFormText

$ws.Range("A29").Value = "IPAIAM070"
$ws.Range("A30").Value = "IPAIAM071"

$ws.Range("D29").Value = "Y"
$ws.Range("D30").Value = "Y"

# Row heights match the source workbook's auto-computed wrap heights.
$ws.Rows.Item(29).RowHeight = 60
$ws.Rows.Item(30).RowHeight = 135

# Move the visible selection down to the newly added data, mirroring the
# scrolled viewport recorded in the edited workbook.
$ws.Range("D29").Select()
